$wb = $excel.ActiveWorkbook

# Source "Patients" sheet and the new "PatientsShifted" sheet (added after
# the last existing sheet, so it lands at the end of the tab strip and
# becomes the active tab - same as in the target workbook).
$src = $wb.Worksheets.Item("Patients")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "PatientsShifted"

# Copy the Patients A1:O5 block into PatientsShifted shifted one column right
# and five rows down (B6:P10), cell by cell so that only cells that actually
# carry data or formatting in the source show up in the destination - this
# matches the sparse layout of the source range (several cells are blank).
for ($r = 1; $r -le 5; $r++) {
  for ($c = 1; $c -le 15; $c++) {
    $srcCell = $src.Cells.Item($r, $c)
    if (($srcCell.Value2 -ne $null) -or ($srcCell.NumberFormat -ne "General")) {
      $dstCell = $new.Cells.Item($r + 5, $c + 1)
      $srcCell.Copy($dstCell)
    }
  }
}

# Match the source column widths (col D -> E, col I -> J).
$new.Columns.Item(5).ColumnWidth = $src.Columns.Item(4).ColumnWidth
$new.Columns.Item(10).ColumnWidth = $src.Columns.Item(9).ColumnWidth

# The new sheet is the active tab, so it should carry the selection marker;
# Excel moves tabSelected off of "Patients" automatically once this sheet
# becomes selected.
[void]$new.Range("B13").Select()
